# Automated map update: the record that used to be on row 79 (Caso 7051,
# "MORENO, JOSE MARIA AV. 345") has been removed from the source feed.
# Deleting the whole row shifts every following record up by one, which
# matches the row-by-row shift seen across rows 79-86 in the new export
# (old row 80 -> new row 79, ..., old row 87 -> new row 86) and shrinks
# the sheet's used range from A1:P87 down to A1:P86.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AYKO")
$ws.Rows.Item(79).Delete()
